$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Item values in A2 and A4 (Aniseed Syrup <-> Chang)
$a2 = $ws.Range("A2").Value()
$a4 = $ws.Range("A4").Value()
$ws.Range("A2").Value = $a4
$ws.Range("A4").Value = $a2

# Update the active selection to C4, matching the new saved view state
$ws.Range("C4").Select()
